$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 93; existing rows 93-207 shift down to 94-208
$ws.Rows(93).Insert()

# Populate the newly inserted row 93 with the new record's data
$ws.Cells.Item(93, 1).Value2 = 5
$ws.Cells.Item(93, 2).Value2 = 'Macroferia Regional de Talca'
$ws.Cells.Item(93, 3).Value2 = 'Maule'
$ws.Cells.Item(93, 4).Value2 = 44483
$ws.Cells.Item(93, 5).Value2 = 7
$ws.Cells.Item(93, 6).Value2 = 100112032
$ws.Cells.Item(93, 7).Value2 = 'Zapallo italiano'
$ws.Cells.Item(93, 8).Value2 = 'Sin especificar'
$ws.Cells.Item(93, 9).Value2 = 'Primera'
$ws.Cells.Item(93, 10).Value2 = 300
$ws.Cells.Item(93, 11).Value2 = 15000
$ws.Cells.Item(93, 12).Value2 = 15000
$ws.Cells.Item(93, 13).Value2 = 15000
$ws.Cells.Item(93, 14).Value2 = '$/caja 60 unidades'
$ws.Cells.Item(93, 15).Value2 = "Región de O'Higgins"
$ws.Cells.Item(93, 16).Value2 = 250
$ws.Cells.Item(93, 17).Value2 = 60
$ws.Cells.Item(93, 18).Value2 = 'Hortaliza'
